$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.914.08"
$ws.Range("E2").Value = "  +7.75%  "

$ws.Range("D3").Value = "3.132.46"
$ws.Range("E3").Value = "  +5.04%  "

$ws.Range("E4").Value = "  +0.21%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "588.31"
$ws.Range("E5").Value = "  +4.06%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "146.51"
$ws.Range("E6").Value = "  +6.90%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.125.34"
$ws.Range("E8").Value = "  +4.95%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.534"
$ws.Range("E9").Value = "  +3.29%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.157"
$ws.Range("E10").Value = "  +18.43%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.71"
$ws.Range("E11").Value = "  +8.32%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.469"
$ws.Range("E12").Value = "  +4.52%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000255"
$ws.Range("E13").Value = "  +11.53%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.83"
$ws.Range("E14").Value = "  +6.74%  "

$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "3.672.65"
$ws.Range("E16").Value = "  +5.78%  "

$ws.Range("D17").Value = "63.916.85"
$ws.Range("E17").Value = "  +7.79%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.149.00"
$ws.Range("E18").Value = "  +5.76%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$ws.Range("E19").Value = "  +0.79%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "471.72"
$ws.Range("E20").Value = "  +8.41%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.18"
$ws.Range("E21").Value = "  +3.64%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.730"
$ws.Range("E22").Value = "  +0.86%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.54"
$ws.Range("E23").Value = "  +7.12%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.30"
$ws.Range("E24").Value = "  +2.38%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "82.35"
$ws.Range("E25").Value = "  +2.92%  "

$ws.Range("E26").Value = "  +0.01%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.63"
$ws.Range("E27").Value = "  +11.81%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.69"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$ws.Range("E30").Value = "  -0.60%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.83"
$ws.Range("E31").Value = "  +11.01%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "27.02"
$ws.Range("E32").Value = "  +5.07%  "

$ws.Range("E33").Value = "  +5.68%  "

$ws.Range("D34").Value = "0.0₃0871"
$ws.Range("E34").Value = "  +14.06%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.41"
$ws.Range("E35").Value = "  +17.76%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.05"
$ws.Range("E36").Value = "  +6.32%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.36"
$ws.Range("E37").Value = "  +19.90%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "6.12"
$ws.Range("E38").Value = "  +3.98%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "50.87"
$ws.Range("E39").Value = "  +4.76%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "443.37"
$ws.Range("E40").Value = "  +12.31%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "8.73"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "2.930.56"
$ws.Range("E42").Value = "  +8.16%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0371"
$ws.Range("E43").Value = "  +5.94%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$ws.Range("E44").Value = "  +8.68%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.282"
$ws.Range("E45").Value = "  +13.13%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$ws.Range("E46").Value = "  +11.17%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "34.79"
$ws.Range("E48").Value = "  +0.95%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "122.88"
$ws.Range("E49").Value = "  +0.09%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$ws.Range("E50").Value = "  +2.20%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "24.75"
$ws.Range("E51").Value = "  +7.14%  "

